# New crime data collected — update the weekly CompStat report
# (63rd Precinct, week of 12/12/2022 - 12/18/2022, Volume 29 Number 50)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text: report volume/number and the date range covered.
# ---------------------------------------------------------------------
$ws.Range("A8").Value2  = "Volume 29   Number  50"
$ws.Range("C9").Value2  = "Report Covering the Week  12/12/2022  Through  12/18/2022"

# ---------------------------------------------------------------------
# Number formats used by the data block (match existing workbook styles)
# ---------------------------------------------------------------------
$fmtCount = "#,##0"
$fmtPct   = "#,##0.0;""-""#,##0.0"
$fmtGeneral = "General"

# ---------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------
$ws.Range("C15").NumberFormat = $fmtCount
$ws.Range("C15").Value = 1
$ws.Range("D15").NumberFormat = $fmtCount
$ws.Range("D15").Value = 1
$ws.Range("E15").NumberFormat = $fmtPct
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 10
$ws.Range("J15").Value = 16
$ws.Range("K15").Value = -37.5
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -41.176470588235
$ws.Range("N15").Value = -58.333333333333

# ---------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 110
$ws.Range("J16").Value = 77
$ws.Range("K16").Value = 42.857142857142
$ws.Range("L16").Value = 15.789473684210
$ws.Range("M16").Value = -49.074074074074
$ws.Range("N16").Value = -83.79970544919

# ---------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -60
$ws.Range("F17").Value = 15
$ws.Range("G17").Value = 11
$ws.Range("H17").Value = 36.363636363636
$ws.Range("I17").Value = 165
$ws.Range("J17").Value = 161
$ws.Range("K17").Value = 2.484472049689
$ws.Range("L17").Value = 12.244897959183
$ws.Range("M17").Value = 27.906976744186
$ws.Range("N17").Value = -46.428571428571

# ---------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------
$ws.Range("C18").Value = 3
$ws.Range("D18").NumberFormat = $fmtCount
$ws.Range("D18").Value = 2
$ws.Range("E18").NumberFormat = $fmtPct
$ws.Range("E18").Value = 50
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = 57.142857142857
$ws.Range("I18").Value = 118
$ws.Range("J18").Value = 99
$ws.Range("K18").Value = 19.191919191919
$ws.Range("L18").Value = 24.210526315789
$ws.Range("M18").Value = -55.471698113207
$ws.Range("N18").Value = -89.65819456617

# ---------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 15
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 51
$ws.Range("G19").Value = 49
$ws.Range("H19").Value = 4.081632653061
$ws.Range("I19").Value = 519
$ws.Range("J19").Value = 399
$ws.Range("K19").Value = 30.075187969924
$ws.Range("L19").Value = 31.392405063291
$ws.Range("M19").Value = 1.169590643274
$ws.Range("N19").Value = -17.619047619047

# ---------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 16
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = 128.571428571429
$ws.Range("I20").Value = 132
$ws.Range("J20").Value = 82
$ws.Range("K20").Value = 60.975609756097
$ws.Range("L20").Value = 23.364485981308
$ws.Range("M20").Value = -15.384615384615
$ws.Range("N20").Value = -95.009451795841

# ---------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 27
$ws.Range("D21").Value = 26
$ws.Range("E21").Value = 3.846153846153
$ws.Range("F21").Value = 103
$ws.Range("G21").Value = 84
$ws.Range("H21").Value = 22.619047619047
$ws.Range("I21").Value = 1055
$ws.Range("J21").Value = 835
$ws.Range("K21").Value = 26.347305389221
$ws.Range("L21").Value = 23.971797884841
$ws.Range("M21").Value = -18.721109399075
$ws.Range("N21").Value = -80.610181951847

# ---------------------------------------------------------------------
# Row 23 - Transit
# ---------------------------------------------------------------------
# "0" looks numeric, so force text format first to stop Excel from
# auto-converting it to a number, matching the original literal text cell.
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "0"
$ws.Range("D23").NumberFormat = $fmtGeneral
$ws.Range("E23").NumberFormat = $fmtGeneral
$ws.Range("E23").Value2 = "***.*"
$ws.Range("F23").Value = 2
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = -33.333333333333
$ws.Range("I23").Value = 31
$ws.Range("J23").Value = 28
$ws.Range("K23").Value = 10.714285714285
$ws.Range("L23").Value = 63.157894736842
$ws.Range("M23").Value = 10.714285714285

# ---------------------------------------------------------------------
# Row 24 - Housing
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 24
$ws.Range("D24").Value = 17
$ws.Range("E24").Value = 41.176470588235
$ws.Range("F24").Value = 93
$ws.Range("G24").Value = 66
$ws.Range("H24").Value = 40.909090909090
$ws.Range("I24").Value = 1111
$ws.Range("J24").Value = 788
$ws.Range("K24").Value = 40.989847715736
$ws.Range("L24").Value = 29.487179487179
$ws.Range("M24").Value = 12.449392712550

# ---------------------------------------------------------------------
# Row 25 - Petit Larceny
# ---------------------------------------------------------------------
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = -40
$ws.Range("F25").Value = 24
$ws.Range("G25").Value = 23
$ws.Range("H25").Value = 4.347826086956
$ws.Range("I25").Value = 235
$ws.Range("J25").Value = 248
$ws.Range("K25").Value = -5.241935483870
$ws.Range("L25").Value = 15.196078431372
$ws.Range("M25").Value = -31.286549707602

# ---------------------------------------------------------------------
# Row 26 - Misd. Assault
# ---------------------------------------------------------------------
$ws.Range("C26").NumberFormat = $fmtCount
$ws.Range("C26").Value = 1
$ws.Range("D26").NumberFormat = $fmtCount
$ws.Range("D26").Value = 1
$ws.Range("E26").NumberFormat = $fmtPct
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 2
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 15
$ws.Range("J26").Value = 26
$ws.Range("K26").Value = -42.307692307692
$ws.Range("L26").Value = 7.142857142857

# ---------------------------------------------------------------------
# Row 27 - UCR Rape*
# ---------------------------------------------------------------------
$ws.Range("C27").NumberFormat = $fmtCount
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 4
$ws.Range("E27").Value = -75
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -57.142857142857
$ws.Range("I27").Value = 32
$ws.Range("J27").Value = 32
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 68.421052631578

# ---------------------------------------------------------------------
# Row 28 - Other Sex Crimes
# ---------------------------------------------------------------------
$ws.Range("M28").Value = -30

# ---------------------------------------------------------------------
# Row 29 - Shooting Vic.
# ---------------------------------------------------------------------
$ws.Range("M29").Value = -41.176470588235
